$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rotate the data rows 2-4: new row2 = old row4, new row3 = old row2, new row4 = old row3
# Capture old values first (use Value2 for reliable scalar reads)
$oldD2 = $ws.Range("D2").Value2
$oldJ2 = $ws.Range("J2").Value2
$oldK2 = $ws.Range("K2").Value2
$oldL2 = $ws.Range("L2").Value2
$oldM2 = $ws.Range("M2").Value2
$oldP2 = $ws.Range("P2").Value2

$oldD3 = $ws.Range("D3").Value2
$oldJ3 = $ws.Range("J3").Value2
$oldK3 = $ws.Range("K3").Value2
$oldL3 = $ws.Range("L3").Value2
$oldM3 = $ws.Range("M3").Value2
$oldP3 = $ws.Range("P3").Value2

$oldD4 = $ws.Range("D4").Value2
$oldJ4 = $ws.Range("J4").Value2
$oldK4 = $ws.Range("K4").Value2
$oldL4 = $ws.Range("L4").Value2
$oldM4 = $ws.Range("M4").Value2
$oldP4 = $ws.Range("P4").Value2

# Row 2 gets old Row 4 values
$ws.Range("D2").Value = $oldD4
$ws.Range("J2").Value = $oldJ4
$ws.Range("K2").Value = $oldK4
$ws.Range("L2").Value = $oldL4
$ws.Range("M2").Value = $oldM4
$ws.Range("P2").Value = $oldP4

# Row 3 gets old Row 2 values
$ws.Range("D3").Value = $oldD2
$ws.Range("J3").Value = $oldJ2
$ws.Range("K3").Value = $oldK2
$ws.Range("L3").Value = $oldL2
$ws.Range("M3").Value = $oldM2
$ws.Range("P3").Value = $oldP2

# Row 4 gets old Row 3 values
$ws.Range("D4").Value = $oldD3
$ws.Range("J4").Value = $oldJ3
$ws.Range("K4").Value = $oldK3
$ws.Range("L4").Value = $oldL3
$ws.Range("M4").Value = $oldM3
$ws.Range("P4").Value = $oldP3
